$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Required Tasks")
$ws1.Range("F4:G7").ClearContents()

$ws2 = $wb.Worksheets.Item("Desirable Tasks")
$ws2.Range("K11:L15").ClearContents()
